# Finished third 3E trial
# Fill in the results (Loss 1-8 / Epoch Min Loss / Train Time) for the
# third trial of the 3-epoch ("3E") run on the "GPT2 - Epochs" sheet,
# in the M column (and the accompanying K/L "3" / "3(E2)" loss columns
# for rows 39-47), mirroring the pattern already used for the other
# trial columns on that sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GPT2 - Epochs")
$ws.Activate()

# Total train time for the new 3(E3) trial (row 38, column M only).
$ws.Range("M38").Formula = "=9418/60"

# Per-row Loss 1..8 and Epoch Min Loss values for columns K, L, M
# (rows 39 through 47).
$data = @(
    @(39, 7240.2841796875,  5629.93310546875,  5413.2197265625),
    @(40, 6797.28076171875, 5617.25537109375,  5403.31005859375),
    @(41, 6199.29248046875, 5586.822265625,    5394.77587890625),
    @(42, 6010.1328125,     5566.54296875,     5379.45458984375),
    @(43, 5892.91845703125, 5522.98388671875,  5361.7548828125),
    @(44, 5823.51904296875, 5507.3515625,      5354.36181640625),
    @(45, 5780.26123046875, 5483.955078125,    5328.583984375),
    @(46, 5720.65576171875, 5457.41259765625,  5315.52001953125),
    @(47, 5683.27197265625, 5456.05224609375,  5310.7060546875)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 11).Value = $row[1]   # column K
    $ws.Cells.Item($r, 12).Value = $row[2]   # column L
    $ws.Cells.Item($r, 13).Value = $row[3]   # column M
}

# Scroll/selection state left behind after entering the new data.
$ws.Range("N40").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 29
$win.ScrollColumn = 3
